$d = $word.ActiveDocument
try {
  $d.Styles(1).Delete()
  Write-Output "deleted"
} catch {
  Write-Output ("ERR: " + $_.Exception.Message)
}
